$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose displayed text changes (Price column D and Volume(1h) column E).
$changedCells = @(
  "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)

# Mark each one as Text format first so Excel keeps the new values as literal
# strings (e.g. "307.73", "2.77%") instead of silently converting the
# numeric-looking text into a number/percentage, matching the original
# inline-string cell contents.
foreach ($addr in $changedCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# New text values taken from the updated symbol list.
$newValues = @{
  'D2' = '307.73'
  'E2' = '2.77%'
  'D3' = '35.92'
  'E3' = '1.80%'
  'D4' = '5.078'
  'E4' = '0.84%'
  'D5' = '0.08111'
  'E5' = '2.36%'
  'D6' = '1.938'
  'E6' = '3.03%'
  'D7' = '4.166'
  'E7' = '3.24%'
  'D8' = '7.834'
  'E8' = '0.62%'
  'D9' = '0.9393'
  'E9' = '1.61%'
  'D10' = '0.1355'
  'E10' = '-7.64%'
  'D11' = '0.1910'
  'E11' = '0.64%'
  'D12' = '0.09260'
  'E12' = '1.29%'
  'D13' = '0.03511'
  'E13' = '1.43%'
  'D14' = '0.09919'
  'E14' = '0.38%'
  'D15' = '0.001452'
  'E15' = '4.42%'
  'D16' = '0.005794'
  'E16' = '0.37%'
  'D17' = '3.609'
  'E17' = '2.92%'
  'D18' = '2.970'
  'E18' = '1.89%'
  'E19' = '1.52%'
  'E20' = '3.95%'
  'D21' = '5.218'
  'E21' = '3.27%'
  'D22' = '0.2530'
  'E22' = '5.21%'
  'D23' = '0.04409'
  'E23' = '-1.27%'
  'D24' = '0.001238'
  'E24' = '1.64%'
  'D25' = '0.004759'
  'E25' = '0.21%'
  'D26' = '0.0001300'
  'E26' = '5.21%'
  'D27' = '0.0003137'
  'E27' = '4.30%'
  'D39' = '0.01999'
  'E39' = '5.65%'
  'D40' = '0.05055'
  'E40' = '7.65%'
  'D41' = '0.01125'
  'E41' = '15.91%'
  'D42' = '0.007612'
  'E42' = '3.56%'
  'D43' = '0.1378'
  'E43' = '4.56%'
  'D44' = '0.002100'
  'E44' = '2.39%'
  'D45' = '0.01131'
  'E45' = '21.14%'
  'D46' = '0.00006340'
  'E46' = '1.15%'
  'D47' = '0.00000000752'
  'E47' = '-0.14%'
  'D48' = '63.57'
  'E48' = '-2.05%'
  'D49' = '0.001194'
  'E49' = '-28.22%'
  'D50' = '0.00002106'
  'E50' = '-0.14%'
  'D51' = '0.0002005'
  'E51' = '-0.14%'
}

foreach ($addr in $changedCells) {
  $ws.Range($addr).Value = $newValues[$addr]
}
